$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.052.13"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.404.60"
$ws.Range("E3").Value = "  +5.20%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.650"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.650"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0941"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.90%  "
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "2.766.61"
$ws.Range("E16").Value = "  +5.22%  "
$ws.Range("D17").Value = "2.397.34"
$ws.Range("E17").Value = "  +4.75%  "
$ws.Range("D18").Value = "43.201.52"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.00%  "
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "276.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0943"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.22%  "
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("E34").Value = "  +7.32%  "
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("E36").Value = "  -3.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("E38").Value = "  -3.26%  "
$ws.Range("E39").Value = "  +3.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.236"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.38%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +44.10%  "
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.508"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.11%  "
$ws.Range("E51").Value = "  +0.95%  "
